{"js": "// Updates the 20x5 arithmetic-problem table, replacing each cell's\n// expression text with the new value from the commit, cell-by-cell in\n// document order (row-major), while preserving each cell's existing\n// paragraph/run formatting (font, size, alignment).\nconst newValues = [\n  [\"45+41=\", \"75+21=\", \"95-57=\", \"56-29=\", \"58+28=\"],\n  [\"11+27=\", \"34+7=\", \"51+11=\", \"90-31=\", \"26+39=\"],\n  [\"34+47=\", \"0+52=\", \"40+17=\", \"21+38=\", \"73-33=\"],\n  [\"71-52=\", \"85-74=\", \"14+79=\", \"58+38=\", \"65-56=\"],\n  [\"72-36=\", \"40+12=\", \"45+52=\", \"8+9=\", \"44-16=\"],\n  [\"72-60=\", \"42+8=\", \"54-13=\", \"88-77=\", \"53+26=\"],\n  [\"98-14=\", \"55+41=\", \"38+3=\", \"96-9=\", \"12+74=\"],\n  [\"54+4=\", \"60-52=\", \"14+45=\", \"2+55=\", \"43+46=\"],\n  [\"64+25=\", \"59-17=\", \"53-28=\", \"48-33=\", \"52-31=\"],\n  [\"64+16=\", \"78-21=\", \"39-4=\", \"16+12=\", \"36+57=\"],\n  [\"77-32=\", \"42+52=\", \"83-52=\", \"90+6=\", \"93-25=\"],\n  [\"19+31=\", \"59+7=\", \"82-26=\", \"35+49=\", \"6+12=\"],\n  [\"26+8=\", \"1+90=\", \"10+48=\", \"9+38=\", \"60-46=\"],\n  [\"56+18=\", \"76-40=\", \"56+2=\", \"25-2=\", \"50-29=\"],\n  [\"24+56=\", \"8+14=\", \"70+10=\", \"31+60=\", \"55-20=\"],\n  [\"76-23=\", \"19-1=\", \"92-87=\", \"13+72=\", \"67-32=\"],\n  [\"80-10=\", \"31-4=\", \"38+61=\", \"56+25=\", \"84-80=\"],\n  [\"13+79=\", \"90-1=\", \"39+27=\", \"54+2=\", \"97-38=\"],\n  [\"20+72=\", \"90-48=\", \"84-27=\", \"70-44=\", \"41+38=\"],\n  [\"97-38=\", \"6+83=\", \"56+33=\", \"26+4=\", \"13+44=\"]\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n// `columnCount` isn't populated by this host; derive the grid shape from\n// `values` instead (a 2D array mirroring the table's rows/cells).\ntable.load(\"rowCount,values\");\nawait context.sync();\n\nconst rowCount = table.rowCount;\nfor (let r = 0; r < rowCount; r++) {\n  const colCount = newValues[r].length;\n  for (let c = 0; c < colCount; c++) {\n    const cell = table.getCell(r, c);\n    cell.value = newValues[r][c];\n  }\n}\n\nawait context.sync();\n", "ps1": "# Updates the 20x5 arithmetic-problem table, replacing each cell's\n# expression text with the new value from the commit, cell-by-cell in\n# document order (row-major), while preserving each cell's existing\n# paragraph/run formatting (font, size, alignment).\n$newValues = @(\n    @(\"45+41=\", \"75+21=\", \"95-57=\", \"56-29=\", \"58+28=\"),\n    @(\"11+27=\", \"34+7=\", \"51+11=\", \"90-31=\", \"26+39=\"),\n    @(\"34+47=\", \"0+52=\", \"40+17=\", \"21+38=\", \"73-33=\"),\n    @(\"71-52=\", \"85-74=\", \"14+79=\", \"58+38=\", \"65-56=\"),\n    @(\"72-36=\", \"40+12=\", \"45+52=\", \"8+9=\", \"44-16=\"),\n    @(\"72-60=\", \"42+8=\", \"54-13=\", \"88-77=\", \"53+26=\"),\n    @(\"98-14=\", \"55+41=\", \"38+3=\", \"96-9=\", \"12+74=\"),\n    @(\"54+4=\", \"60-52=\", \"14+45=\", \"2+55=\", \"43+46=\"),\n    @(\"64+25=\", \"59-17=\", \"53-28=\", \"48-33=\", \"52-31=\"),\n    @(\"64+16=\", \"78-21=\", \"39-4=\", \"16+12=\", \"36+57=\"),\n    @(\"77-32=\", \"42+52=\", \"83-52=\", \"90+6=\", \"93-25=\"),\n    @(\"19+31=\", \"59+7=\", \"82-26=\", \"35+49=\", \"6+12=\"),\n    @(\"26+8=\", \"1+90=\", \"10+48=\", \"9+38=\", \"60-46=\"),\n    @(\"56+18=\", \"76-40=\", \"56+2=\", \"25-2=\", \"50-29=\"),\n    @(\"24+56=\", \"8+14=\", \"70+10=\", \"31+60=\", \"55-20=\"),\n    @(\"76-23=\", \"19-1=\", \"92-87=\", \"13+72=\", \"67-32=\"),\n    @(\"80-10=\", \"31-4=\", \"38+61=\", \"56+25=\", \"84-80=\"),\n    @(\"13+79=\", \"90-1=\", \"39+27=\", \"54+2=\", \"97-38=\"),\n    @(\"20+72=\", \"90-48=\", \"84-27=\", \"70-44=\", \"41+38=\"),\n    @(\"97-38=\", \"6+83=\", \"56+33=\", \"26+4=\", \"13+44=\")\n)\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n$rowCount = $t.Rows.Count\n\nfor ($r = 1; $r -le $rowCount; $r++) {\n    $rowValues = $newValues[$r - 1]\n    $colCount = $rowValues.Length\n    for ($c = 1; $c -le $colCount; $c++) {\n        $cell = $t.Cell($r, $c)\n        $cell.Range.Text = $rowValues[$c - 1]\n    }\n}\n"}
